$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point drift on the existing A2 timestamp value
$ws.Range("A2").Value = 45878.04184609954

# Append the new data row (row 3)
$ws.Range("A3").Value = 45878.08351149064
$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 37
$ws.Range("D3").Value = 13.33
$ws.Range("E3").Value = 91.53
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2.69
$ws.Range("H3").Value = "E"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "02:00:15"

# Match the style of the previous date cell (A2) so formatting stays consistent
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
